$wb = $excel.ActiveWorkbook

# Rename the existing sheet to pass1
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "pass1"

# Add a new sheet "pass2" after pass1
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "pass2"

# Fill in the new KO values
$ws2.Range("A1").Value = "K01133"
$ws2.Range("A2").Value = "K05846"
$ws2.Range("A3").Value = "K05847"

# Match formatting style used by the later rows in pass1 (style index 1 / fontId 1)
$ws1.Range("A23").Copy()
$ws2.Range("A1:A3").PasteSpecial(-4122)

# Make pass2 the active (selected) sheet/tab
$ws2.Activate()
$ws2.Range("A1").Select()
